$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$dede = $wb.Worksheets.Item("de-de")
$overview = $wb.Worksheets.Item("Overview")

$nineDfUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d827816f97303e62d39ca2926ba948eb15741d03/e2e/9df1ad5d-a0da-4e78-bbec-6cc47ee6e311.md"
$e40Url = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/d827816f97303e62d39ca2926ba948eb15741d03/e2e/e40dd20d-ce14-4908-879f-0d3f7224bf66.md"

# Overview sheet: status strings updated for zh-cn/de-de columns
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"

# Widen the zh-cn/de-de status columns on Overview
$overview.Columns.Item(5).ColumnWidth = 29.9777047293527
$overview.Columns.Item(6).ColumnWidth = 29.9777047293527

# zh-cn sheet: widen Status / Latest Target File / Latest Handback File columns
$zhcn.Columns.Item(3).ColumnWidth = 29.9777047293527
$zhcn.Columns.Item(9).ColumnWidth = 40
$zhcn.Columns.Item(10).ColumnWidth = 40

# de-de sheet: widen Status / Latest Target File / Latest Handback File columns
$dede.Columns.Item(3).ColumnWidth = 29.9777047293527
$dede.Columns.Item(9).ColumnWidth = 40
$dede.Columns.Item(10).ColumnWidth = 40

# zh-cn rows: fill in Latest Target File / Handback File / Handback DateTime
$zhcn.Range("I2").Value = "9df1ad5d-a0da-4e78-bbec-6cc47ee6e311.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $nineDfUrl, "", "", "9df1ad5d-a0da-4e78-bbec-6cc47ee6e311.md")
$zhcn.Range("J2").Value = "9df1ad5d-a0da-4e78-bbec-6cc47ee6e311.4464ee521c665dde2c76c1fac8c000ca9bdd849e.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-16 00:22:40"

$zhcn.Range("I3").Value = "e40dd20d-ce14-4908-879f-0d3f7224bf66.md"
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $e40Url, "", "", "e40dd20d-ce14-4908-879f-0d3f7224bf66.md")
$zhcn.Range("J3").Value = "e40dd20d-ce14-4908-879f-0d3f7224bf66.5003188178e531960ec2b0a1a5a82695514fc486.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-16 00:22:40"

# de-de rows: fill in Latest Target File / Handback File / Handback DateTime
$dede.Range("I2").Value = "9df1ad5d-a0da-4e78-bbec-6cc47ee6e311.md"
$dede.Hyperlinks.Add($dede.Range("I2"), $nineDfUrl, "", "", "9df1ad5d-a0da-4e78-bbec-6cc47ee6e311.md")
$dede.Range("J2").Value = "9df1ad5d-a0da-4e78-bbec-6cc47ee6e311.4464ee521c665dde2c76c1fac8c000ca9bdd849e.de-de.xlf"
$dede.Range("K2").Value = "2016-08-16 00:22:47"

$dede.Range("I3").Value = "e40dd20d-ce14-4908-879f-0d3f7224bf66.md"
$dede.Hyperlinks.Add($dede.Range("I3"), $e40Url, "", "", "e40dd20d-ce14-4908-879f-0d3f7224bf66.md")
$dede.Range("J3").Value = "e40dd20d-ce14-4908-879f-0d3f7224bf66.5003188178e531960ec2b0a1a5a82695514fc486.de-de.xlf"
$dede.Range("K3").Value = "2016-08-16 00:22:47"
